$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the formatting used for the existing date cell (C2) so the new
# date cells can reuse the same number format (built-in short date, s="1").
$ws.Range("C2").Copy()
$ws.Range("C3:C9").PasteSpecial(-4122)  # xlPasteFormats

# Row 2 - groceries
$ws.Range("A2").Value = "groceries"
$ws.Range("B2").Value = 2000
$ws.Range("C2").Value = 45828.22928240741

# Row 3 - food
$ws.Range("A3").Value = "food"
$ws.Range("B3").Value = 200
$ws.Range("C3").Value = 45828.22928240741

# Row 4 - investment
$ws.Range("A4").Value = "investment"
$ws.Range("B4").Value = 2000
$ws.Range("C4").Value = 45827.22928240741

# Row 5 - food
$ws.Range("A5").Value = "food"
$ws.Range("B5").Value = 200
$ws.Range("C5").Value = 45819.72928240741

# Row 6 - petrol
$ws.Range("A6").Value = "petrol"
$ws.Range("B6").Value = 500
$ws.Range("C6").Value = 45818.72928240741

# Row 7 - RENT
$ws.Range("A7").Value = "RENT"
$ws.Range("B7").Value = 10000
$ws.Range("C7").Value = 45809.22928240741

# Row 8 - petrol
$ws.Range("A8").Value = "petrol"
$ws.Range("B8").Value = 500
$ws.Range("C8").Value = 45787.72928240741

# Row 9 - bill (previously row 2, amount changed from 1234565 to 4400)
$ws.Range("A9").Value = "bill"
$ws.Range("B9").Value = 4400
$ws.Range("C9").Value = 45785.72928240741
